# This script updates the cryptocurrency price/volume table with freshly
# scraped values (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) format, used to restore the
# original style on price cells after temporarily marking them as Text so
# Excel does not reinterpret numeric-looking strings (e.g. "244.40",
# "0.5080", "66.00") as numbers and silently drop significant trailing
# zeros / precision.
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "29.579.28"
$ws.Range("E2").Value = "  +1.86%  "
$ws.Range("D3").Value = "1.843.69"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9983"
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").Value = "  -0.75%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.40"
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6312"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  +2.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9994"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("E8").Value = "  +1.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2943"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  +2.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.87"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  +4.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07675"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "1.842.78"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.025"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("E14").Value = "  +3.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "84.16"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  +3.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009324"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  +4.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.976"
$ws.Range("D17").Style = $defaultStyle
$ws.Range("E17").Value = "  +2.37%  "
$ws.Range("D18").Value = "29.535.70"
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("D19").Value = "2.089.86"
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.03"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.58"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9993"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.362"
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = "  +3.92%  "
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.70"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1419"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.522"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.497"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E30").Value = "  +8.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.254"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  +3.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.145"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +1.67%  "
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("E34").Value = "  +2.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.149"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7262"
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.609"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.879"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  +2.31%  "
$ws.Range("D39").Value = "1.222.07"
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01769"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.300"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9185"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +4.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  -0.52%  "
$ws.Range("D44").Value = "2.009.01"
$ws.Range("E44").Value = "  +2.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.91"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.00"
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = "  +2.89%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000121"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5080"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.261"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4072"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1134"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  +3.95%  "

Write-Host "Applied all crypto list updates"
